$wb = $excel.ActiveWorkbook

# --- Sheet1: Postulantes ---
$ws1 = $wb.Worksheets.Item("Postulantes")

# Column A (numbers)
$ws1.Cells.Item(2, 1).Value = 1
$ws1.Cells.Item(3, 1).Value = 2
$ws1.Cells.Item(4, 1).Value = 3
$ws1.Cells.Item(5, 1).Value = 4
$ws1.Cells.Item(6, 1).Value = 5
$ws1.Cells.Item(7, 1).Value = 6
$ws1.Cells.Item(8, 1).Value = 7
$ws1.Cells.Item(9, 1).Value = 8
$ws1.Cells.Item(10, 1).Value = 9
$ws1.Cells.Item(11, 1).Value = 10
$ws1.Cells.Item(12, 1).Value = 11

# Column B (text)
$ws1.Cells.Item(2, 2).Value = "ASDasd"
$ws1.Cells.Item(3, 2).Value = "dassda"
$ws1.Cells.Item(4, 2).Value = "dassda"
$ws1.Cells.Item(5, 2).Value = "ASDasd"
$ws1.Cells.Item(6, 2).Value = "dassda"
$ws1.Cells.Item(7, 2).Value = "dassda"
$ws1.Cells.Item(8, 2).Value = "ASDasd"
$ws1.Cells.Item(9, 2).Value = "dassda"
$ws1.Cells.Item(10, 2).Value = "dassda"
$ws1.Cells.Item(11, 2).Value = "ASDasd"
$ws1.Cells.Item(12, 2).Value = "dassda"

# Column C (text)
$ws1.Cells.Item(3, 3).Value = "sad"
$ws1.Cells.Item(4, 3).Value = "a"
$ws1.Cells.Item(8, 3).Value = "sadadsa"
$ws1.Cells.Item(5, 3).Value = "asd"
$ws1.Cells.Item(7, 3).Value = "dsaads"
$ws1.Cells.Item(6, 3).Value = "sda"
$ws1.Cells.Item(9, 3).Value = "dsadas"
$ws1.Cells.Item(10, 3).Value = "dadsads"
$ws1.Cells.Item(11, 3).Value = "adsasd"
$ws1.Cells.Item(12, 3).Value = "adssda"

# Column D (Grupo Ocupacional) - all CHOFER
$ws1.Cells.Item(2, 4).Value = "CHOFER"
$ws1.Cells.Item(3, 4).Value = "CHOFER"
$ws1.Cells.Item(4, 4).Value = "CHOFER"
$ws1.Cells.Item(5, 4).Value = "CHOFER"
$ws1.Cells.Item(6, 4).Value = "CHOFER"
$ws1.Cells.Item(7, 4).Value = "CHOFER"
$ws1.Cells.Item(8, 4).Value = "CHOFER"
$ws1.Cells.Item(9, 4).Value = "CHOFER"
$ws1.Cells.Item(10, 4).Value = "CHOFER"
$ws1.Cells.Item(11, 4).Value = "CHOFER"
$ws1.Cells.Item(12, 4).Value = "CHOFER"

# --- Sheet2: Plazas ---
$ws2 = $wb.Worksheets.Item("Plazas")

# Column A
$ws2.Cells.Item(2, 1).Value = "GOF"
$ws2.Cells.Item(3, 1).Value = "GOF"
$ws2.Cells.Item(4, 1).Value = "CENATE"
$ws2.Cells.Item(5, 1).Value = "CENATE"
$ws2.Cells.Item(6, 1).Value = "CENATE"

# Column B
$ws2.Cells.Item(2, 2).Value = "olivos"
$ws2.Cells.Item(3, 2).Value = "lince"
$ws2.Cells.Item(4, 2).Value = "olivos"
$ws2.Cells.Item(5, 2).Value = "smp"
$ws2.Cells.Item(6, 2).Value = "centro"

# Column D
$ws2.Cells.Item(2, 4).Value = 5
$ws2.Cells.Item(3, 4).Value = 12
$ws2.Cells.Item(4, 4).Value = 2
$ws2.Cells.Item(5, 4).Value = 1
$ws2.Cells.Item(6, 4).Value = 6

# Sheet1 row 2 col C is entered last, matching the original authoring sequence
$ws1.Cells.Item(2, 3).Value = "aaaaa"

# Page setup for Postulantes sheet (now printable as portrait)
$ws1.PageSetup.Orientation = 1

# Final selections per sheet, then leave Postulantes as the active sheet/tab
$ws2.Range("G7").Select()
$ws1.Select()
$ws1.Range("A13").Select()

$wb.Save()
